$d = $word.ActiveDocument

# Locate the "Docente(s) Responsável(eis)" heading paragraph so the new
# professor list can be inserted right after it.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Docente(s) Responsável(eis)*") {
        $target = $p
        break
    }
}
if ($target -eq $null) {
    throw "Could not find the 'Docente(s) Responsável(eis)' paragraph"
}

# Create a new, empty paragraph right after the heading.
$null = $target.Range.InsertParagraphAfter()
$newPara = $target.Next()

# Fill the new paragraph via a raw OOXML fragment so the two professor
# names land in separate runs (first run carries the line break), matching
# how the rest of the document structures multi-line list entries.
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
       '<pkg:xmlData>' +
       '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:body>' +
       '<w:p>' +
       '<w:pPr><w:pStyle w:val="ListBullet"/></w:pPr>' +
       '<w:r><w:t>5817344 - Livia Melo Carneiro</w:t><w:br/></w:r>' +
       '<w:r><w:t>6310296 - Patrícia Caroline Molgero Da Rós</w:t></w:r>' +
       '</w:p>' +
       '</w:body></w:document>' +
       '</pkg:xmlData></pkg:part></pkg:package>'

$newPara.Range.InsertXML($xml)
